$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.494.54'
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = '2.545.43'
$ws.Range("E3").Value = '  +4.72%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.99'
$ws.Range("E5").Value = '  +2.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.36'
$ws.Range("E6").Value = '  +9.12%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").Value = '2.544.05'
$ws.Range("E9").Value = '  +4.69%  '
$ws.Range("E10").Value = '  +2.68%  '
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("E13").Value = '  +3.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.38'
$ws.Range("E14").Value = '  +8.75%  '
$ws.Range("D15").Value = '3.001.38'
$ws.Range("E15").Value = '  +4.75%  '
$ws.Range("D16").Value = '63.422.44'
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").Value = '2.503.53'
$ws.Range("E18").Value = '  +2.61%  '
$ws.Range("E19").Value = '  +4.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '340.75'
$ws.Range("E20").Value = '  -1.62%  '
$ws.Range("E21").Value = '  +4.45%  '
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.23'
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("E26").Value = '  +5.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.51'
$ws.Range("E27").Value = '  +13.21%  '
$ws.Range("E28").Value = '  +3.99%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.12'
$ws.Range("E30").Value = '  +11.91%  '
$ws.Range("D31").Value = '0.0₃0836'
$ws.Range("E31").Value = '  +6.57%  '
$ws.Range("E32").Value = '  +3.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '178.15'
$ws.Range("E33").Value = '  +3.80%  '
$ws.Range("E34").Value = '  +9.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '419.27'
$ws.Range("E35").Value = '  +14.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.406'
$ws.Range("E36").Value = '  +2.70%  '
$ws.Range("E37").Value = '  +3.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.47'
$ws.Range("E38").Value = '  +0.22%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("E40").Value = '  +4.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.01'
$ws.Range("E42").Value = '  +2.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '154.86'
$ws.Range("E43").Value = '  +6.10%  '
$ws.Range("E44").Value = '  +4.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.16'
$ws.Range("E45").Value = '  +3.09%  '
$ws.Range("E46").Value = '  +4.08%  '
$ws.Range("E47").Value = '  +3.00%  '
$ws.Range("E48").Value = '  +9.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0969'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.64'
$ws.Range("E50").Value = '  +4.45%  '
$ws.Range("E51").Value = '  +7.44%  '
